$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update subject names (cells keep their original quote-prefixed text style,
# so re-enter the values with a leading apostrophe just like the original
# authoring did, to avoid Excel dropping the quote-prefix cell style):
# - row 4 (B4): "KH Xa Hoi" -> "KHXH"
# - row 5 (B5): "Tieng Anh" (unchanged text, but shared-string slot moves)
# - row 6 (B6): "KH Tu Nhien" -> "KHTN"
# - row 7 (B7): "Tin Hoc" (unchanged text, but shared-string slot moves)
$ws.Range("B4").Value = "'KHXH"
$ws.Range("B5").Value = "'Tieng Anh"
$ws.Range("B6").Value = "'KHTN"
$ws.Range("B7").Value = "'Tin Hoc"

# Update the sheet view: scroll/zoom + new selection
$ws.Application.ActiveWindow.Zoom = 139
$ws.Application.Goto($ws.Range("A2"), $true)
$ws.Range("C7").Select()
